# Apply the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Figure1C")
$ws3 = $wb.Worksheets.Item("Figure1F")
$ws4 = $wb.Worksheets.Item("Figure2A")
$ws5 = $wb.Worksheets.Item("Figure2B")

# --- Figure1C (sheet 1): fill in previously-empty B/C cells for rows 5-8 ---
$ws1.Range("B5").Value = 0.101136
$ws1.Range("C5").Value = 0.111634
$ws1.Range("B6").Value = 0.22692999999999999
$ws1.Range("C6").Value = 0.06544
$ws1.Range("C7").Value = 0.111634
$ws1.Range("C8").Value = 0.160318

# --- Figure1F (sheet 3): change 1E-4 placeholders to 1E-5 ---
$ws3.Range("B2").Value = 0.00001
$ws3.Range("C2").Value = 0.00001
$ws3.Range("C3").Value = 0.00001
$ws3.Range("C4").Value = 0.00001
$ws3.Range("C5").Value = 0.00001
$ws3.Range("B6").Value = 0.00001
$ws3.Range("C6").Value = 0.00001
$ws3.Range("C7").Value = 0.00001
$ws3.Range("C8").Value = 0.00001

# --- Figure2A (sheet 4): change 1E-4 placeholders to 1E-5 ---
$ws4.Range("C2").Value = 0.00001
$ws4.Range("B3").Value = 0.00001
$ws4.Range("C3").Value = 0.00001
$ws4.Range("C4").Value = 0.00001
$ws4.Range("C5").Value = 0.00001
$ws4.Range("C6").Value = 0.00001
$ws4.Range("C7").Value = 0.00001
$ws4.Range("B8").Value = 0.00001
$ws4.Range("C8").Value = 0.00001
$ws4.Range("C10").Value = 0.00001
$ws4.Range("C11").Value = 0.00001
$ws4.Range("C12").Value = 0.00001
$ws4.Range("C13").Value = 0.00001

# --- Update each sheet's stored selection (cursor position) ---
# Order matters: the last sheet selected becomes the active tab, and
# Figure2A (sheet 4) must remain the active tab, as in the original file.
$ws1.Range("C8").Select() | Out-Null
$ws3.Range("F8").Select() | Out-Null
$ws5.Range("C2").Select() | Out-Null
$ws4.Range("C9").Select() | Out-Null
